# LegacyDatabase/Users.xlsx — rename the sole "data" sheet to "Users" and
# move the saved selection from A5 to B52 (cell A1:G201 data left untouched).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("data")
$ws.Name = "Users"

$ws.Activate()
$ws.Range("B52").Select()
